$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Structural changes
# -----------------------------------------------------------------
# Remove the old "GL Code", "Location", "Title" columns (C:E)
$ws.Range("C1:E1").EntireColumn.Delete()

# Remove the two trailing rows (old rows 6 and 7) that no longer exist
$ws.Range("A6:A7").EntireRow.Delete()

# Insert a new column before the old "Plan" column to hold "Division"
$ws.Range("C1").EntireColumn.Insert()

# -----------------------------------------------------------------
# Header row (row 1)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Notes"
$ws.Range("B1").Value = "Employee Name"
$ws.Range("C1").Value = "Division"
$ws.Range("D1").Value = "Plan"
$ws.Range("E1").Value = "Tier"
$ws.Range("F1").Value = "Funding Amount"

# -----------------------------------------------------------------
# Data rows
#
# The "Funding Amount" column historically stored its numbers as text
# (e.g. the old values "24.62"/"123.10" were plain strings, not
# numbers), so force that column to text formatting before writing the
# new amounts to keep them stored the same way.
# -----------------------------------------------------------------
$ws.Range("F2:F5").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Joined"
$ws.Range("B2").Value = "John Jones"
$ws.Range("C2").Value = "HR"
$ws.Range("D2").Value = "Carrier1"
$ws.Range("E2").Value = "Tier3"
$ws.Range("F2").Value = "10418.43"

# Row 3
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "Michelle Johnson"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Carrier2"
$ws.Range("E3").Value = "Tier1"
$ws.Range("F3").Value = "17.63"

# Row 4
$ws.Range("A4").Value = "Joined"
$ws.Range("B4").Value = "Jane Doe"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Carrier1"
$ws.Range("E4").Value = "Tier4"
$ws.Range("F4").Value = "406.38"

# Row 5 - only the Funding Amount remains populated
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "10842.44"
